# ------------------------------------------------------------------
# C5-PowerPoint.pptx edit
#
#  1) The table on slide 6 gets a new (built-in) table style applied
#     instead of the bespoke "Table_0" style that ships in
#     ppt/tableStyles.xml.
#
#  2) The colour scheme that is baked into ppt/theme/theme1.xml (the
#     theme used by the slide master / slides, originally the
#     "Integral" palette) is swapped for the colour scheme that used
#     to live in ppt/theme/theme2.xml (the stock "Office Theme"
#     palette used by the notes master). The font scheme / format
#     scheme are identical between the two themes, so re-pointing the
#     twelve theme colour slots reproduces the content swap described
#     by the diff.
# ------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 6 ------------------------------------
for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{30492B1E-891A-4AE0-B700-CC8930C6D87E}")
        }
    }
}

# --- 2. Swap the slide-master theme colours for the former notes- -------
#        master ("Office Theme") palette.
$slideForTheme = $p.Slides.Item(1)
$themeColors = $slideForTheme.ThemeColorScheme

function ColorFromHex($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Slot order matches a:clrScheme: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeHex = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ColorFromHex($officeThemeHex[$i - 1])
}
